$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full (column A / C) names for the 24 new rows (280-303), in row order.
$names = @(
    "Bristol City",
    "Coventry City",
    "Hull City",
    "Luton Town",
    "Preston North End",
    "West Bromwich Albion",
    "Birmingham City",
    "Burnley",
    "Cardiff City",
    "Huddersfield Town",
    "Norwich City",
    "Wigan Athletic",
    "Reading",
    "Blackburn Rovers",
    "Blackpool",
    "Middlesbrough",
    "Millwall",
    "Queens Park Rangers",
    "Rotherham United",
    "Sheffield United",
    "Stoke City",
    "Sunderland",
    "Swansea City",
    "Watford"
)

# FBRef short-form names for column B, in row order (same index as $names).
$fbref = @(
    "Bristol City",
    "Coventry City",
    "Hull City",
    "Luton Town",
    "Preston",
    "West Brom",
    "Birmingham City",
    "Burnley",
    "Cardiff City",
    "Huddersfield",
    "Norwich City",
    "Wigan Athletic",
    "Reading",
    "Blackburn",
    "Blackpool",
    "Middlesbrough",
    "Millwall",
    "QPR",
    "Rotherham Utd",
    "Sheffield Utd",
    "Stoke City",
    "Sunderland",
    "Swansea City",
    "Watford"
)

$startRow = 280

# Fill column A (Name) first, top to bottom.
for ($i = 0; $i -lt $names.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $names[$i]
}

# Fill column C (Bovada) next, top to bottom (mirrors column A).
for ($i = 0; $i -lt $names.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 3).Value = $names[$i]
}

# Fill column B (FBRef): the handful of rows whose short name differs from
# the full name were entered first (in this particular order), then the
# remaining rows (whose FBRef name equals the full name) were filled in.
$specialRows = @(289, 298, 293, 284, 297, 299, 285)
foreach ($r in $specialRows) {
    $idx = $r - $startRow
    $ws.Cells.Item($r, 2).Value = $fbref[$idx]
}

for ($i = 0; $i -lt $names.Count; $i++) {
    $r = $startRow + $i
    if ($specialRows -notcontains $r) {
        $ws.Cells.Item($r, 2).Value = $fbref[$i]
    }
}

# Match the saved view state: the last cell of the newly-added data selected.
$ws.Range("A303").Select()
